# Actualización automática 2025-11-25 10:30:08
$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M32").Value = 1946.43
$ws1.Range("H49").Value = 916.2
$ws1.Range("I49").Value = 183.4
$ws1.Range("L49").Value = 1334.66
$ws1.Range("M49").Value = 1481.93

# --- Hoja "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F32").Value = 1946.43
$ws2.Range("F49").Value = 3916.19
$ws2.Range("F61").Value = 51052.21

# --- Hoja "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D6").Value = 3203.12
$ws3.Range("E6").Value = -1598.12
$ws3.Range("F6").Value = 1.995713395638629

$ws3.Range("D7").Value = 641.79
$ws3.Range("E7").Value = 244.921016287574
$ws3.Range("F7").Value = 0.7237871056198287

$ws3.Range("D11").Value = 7785.52
$ws3.Range("E11").Value = 8362.48
$ws3.Range("F11").Value = 0.4821352489472381

$ws3.Range("D12").Value = 17683.63
$ws3.Range("E12").Value = 32623.37
$ws3.Range("F12").Value = 0.3515143021845867

$ws3.Range("D14").Value = 54018.15000000001
$ws3.Range("E14").Value = 43843.73766749098
$ws3.Range("F14").Value = 0.5519835278830867
